$d = $word.ActiveDocument
$target = $null
$idx = 0
$targetIdx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*let result = square(number);*") {
        $target = $p
        $targetIdx = $idx
    }
}
if ($target -eq $null) {
    throw "Target paragraph not found"
}
$target.Range.InsertParagraphAfter()
$newp = $d.Paragraphs($targetIdx + 1)
$xml = '<w:p/><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t>Arrow Functions:</w:t></w:r></w:p><w:p><w:r><w:t>More readable and more structured</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Anonymous </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>functions(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>lambda functions)</w:t></w:r></w:p><w:p><w:r><w:t>Without function name but they are assigned to a variable</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Let c</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)=&gt;{</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Console.log(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>“Hi Ashish”)</w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>C(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>o/p:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Hi Ashish</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Arrow function with parameters</w:t></w:r></w:p><w:p><w:r><w:t>let x=(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>x,y</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>)=&gt;{</w:t></w:r></w:p><w:p><w:r><w:t>            console.log(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>x+y</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>x(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>20,30)</w:t></w:r></w:p><w:p><w:r><w:t>o/p:</w:t></w:r></w:p><w:p><w:r><w:t>50</w:t></w:r></w:p><w:p><w:r><w:t>Arrow function with one argument</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Without using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>paranthesis</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>Let greet=x=&gt;console.log(x);</w:t></w:r></w:p><w:p><w:r><w:t>greet(‘hello’)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Arrow function </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>with out</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> argument</w:t></w:r></w:p><w:p><w:r><w:t>let greet=</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>()=</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>&gt;console.log(“hello”)</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>greet(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">//Arrow functions as </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> expression</w:t></w:r></w:p><w:p><w:r><w:t>        let age=18;</w:t></w:r></w:p><w:p><w:r><w:t>        let welcome=(age&lt;18)?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>()=</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>&gt;console.log("hi Gowtham"):</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>()=</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">&gt;console.log("hi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>chran</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>");</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>welcome(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>o/p:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">hi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>chran</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:r><w:t> //synchronous and asynchronous</w:t></w:r></w:p><w:p><w:r><w:t>        function task1</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>(){</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:t>            console.log("Task</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>1:Start</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>")</w:t></w:r></w:p><w:p><w:r><w:t>        }</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>        function task2</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>(){</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:t>            console.log("Task</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>2:Start</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>")</w:t></w:r></w:p><w:p><w:r><w:t>        }</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>        function task3</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>(){</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:t>            console.log("</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Task:Start</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>")</w:t></w:r></w:p><w:p><w:r><w:t>        }</w:t></w:r></w:p><w:p><w:r><w:t>        task1()</w:t></w:r></w:p><w:p><w:r><w:t>        task2()</w:t></w:r></w:p><w:p><w:r><w:t>        task3()</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t>//</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>asynchronomous</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>        //multiple tasks at a time</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>console.log(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">"hi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ashish</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>")</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>setTimeout</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>()=&gt;{</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>console.log(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">"Hi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>iam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> waiting")</w:t></w:r></w:p><w:p><w:r><w:t>        },5000);</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>console.log(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>"I am in lab")</w:t></w:r></w:p><w:p/><w:p><w:r><w:t> //</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javascript</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> callbacks</w:t></w:r></w:p><w:p><w:r><w:t>        //A call back function is a function that is passed argument to another function;</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">        function </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>greet(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>name, callback) {</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>console.log(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>`Hello, ${name}!`);</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>callback(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:r><w:t>        //this calls the callback functions after greeting</w:t></w:r></w:p><w:p><w:r><w:t>        }</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        function </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>sayGoodbye</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>){</w:t></w:r></w:p><w:p><w:r><w:t>            console.log("Goodbye")</w:t></w:r></w:p><w:p><w:r><w:t>        }</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        //call greet and pass </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sayGoodbyeas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the callback</w:t></w:r></w:p><w:p><w:r><w:t>        greet("</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramStart"/><w:r><w:t>",</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sayGoodbye</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>     </w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t>//promise</w:t></w:r></w:p><w:p><w:r><w:t>        //eventual completion</w:t></w:r></w:p><w:p><w:r><w:t>        //Three states</w:t></w:r></w:p><w:p><w:r><w:t>        //</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>1)</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pending</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>:still</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> not completed --continuing</w:t></w:r></w:p><w:p><w:r><w:t>        //</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>2)</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fulfilled</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>:task</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> completed</w:t></w:r></w:p><w:p><w:r><w:t>        //</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>3)</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rejected</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>:The</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> operation is failed</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">        let </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mypromise</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=new Promise((</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>resolve,reject</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>)=&gt;{</w:t></w:r></w:p><w:p><w:r><w:t>            let success=true;</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>            if(success</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>){</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">                </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>resolve(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>"task completed")</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>}else</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>{</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>reject(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>"task failed")</w:t></w:r></w:p><w:p><w:r><w:t>            }</w:t></w:r></w:p><w:p><w:r><w:t>        })</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/>'
$newp.Range.InsertXML($xml)
Write-Output "Inserted block after target paragraph (index $targetIdx)"
